# Append 4 new match rows (rows 132-135) to Sheet1, mirroring the existing
# layout/formatting of the preceding data rows (e.g. row 131).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 131

$rows = @(
    @{ Indice=131; E=45255.6875; F="FC Voluntari";    G=4; H="Farul Constanta"; I=2;
       J=3.3;  K="23/11/2023 11:42"; L=3.14; M="25/11/2023 16:26";
       N=3.26; O="23/11/2023 11:42"; P=3.23; Q="25/11/2023 16:26";
       R=2.25; S="23/11/2023 11:42"; T=2.4;  U="25/11/2023 16:26";
       V="https://www.betexplorer.com/football/romania/liga-1/voluntari-farul-constanta/r3bsAz1G/" }

    @{ Indice=132; E=45255.8125; F="CFR Cluj";         G=0; H="UTA Arad";         I=0;
       J=1.5;  K="23/11/2023 11:42"; L=1.43; M="25/11/2023 19:28";
       N=4.22; O="23/11/2023 11:42"; P=4.68; Q="25/11/2023 19:28";
       R=6.42; S="23/11/2023 11:42"; T=7.35; U="25/11/2023 19:28";
       V="https://www.betexplorer.com/football/romania/liga-1/cfr-cluj-fc-uta-arad/nHEvjF9p/" }

    @{ Indice=133; E=45256.6875; F="Otelul";            G=1; H="Univ. Craiova";    I=3;
       J=3.42; K="23/11/2023 11:42"; L=3.9;  M="26/11/2023 16:22";
       N=3;    O="23/11/2023 11:42"; P=3.08; Q="26/11/2023 16:22";
       R=2.34; S="23/11/2023 11:42"; T=2.15; U="26/11/2023 16:22";
       V="https://www.betexplorer.com/football/romania/liga-1/otelul-univ-craiova/0KIrkZOj/" }

    @{ Indice=134; E=45256.8125; F="Din. Bucuresti";    G=0; H="FCSB";             I=1;
       J=6.59; K="23/11/2023 11:42"; L=6.22; M="26/11/2023 19:21";
       N=4.6;  O="23/11/2023 11:42"; P=4.24; Q="26/11/2023 19:21";
       R=1.45; S="23/11/2023 11:42"; T=1.53; U="26/11/2023 19:21";
       V="https://www.betexplorer.com/football/romania/liga-1/din-bucuresti-fcsb/EN6QDI1i/" }
)

$targetRow = $lastExistingRow
foreach ($row in $rows) {
    $targetRow = $targetRow + 1

    # Copy full cell formatting (style/number-format/font/border) from the
    # corresponding cell one row above so the new row matches the sheet's
    # existing per-column styling (bold+border index column, datetime
    # format on the match-date column, plain default style elsewhere).
    $ws.Range($ws.Cells.Item($lastExistingRow, 1), $ws.Cells.Item($lastExistingRow, 22)).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($targetRow, 1), $ws.Cells.Item($targetRow, 22)).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($targetRow, 1).Value2  = $row.Indice
    $ws.Cells.Item($targetRow, 2).Value2  = "romania"
    $ws.Cells.Item($targetRow, 3).Value2  = "liga-1"
    $ws.Cells.Item($targetRow, 4).Value2  = "2023-2024"
    $ws.Cells.Item($targetRow, 5).Value2  = $row.E
    $ws.Cells.Item($targetRow, 6).Value2  = $row.F
    $ws.Cells.Item($targetRow, 7).Value2  = $row.G
    $ws.Cells.Item($targetRow, 8).Value2  = $row.H
    $ws.Cells.Item($targetRow, 9).Value2  = $row.I
    $ws.Cells.Item($targetRow, 10).Value2 = $row.J
    $ws.Cells.Item($targetRow, 11).Value2 = $row.K
    $ws.Cells.Item($targetRow, 12).Value2 = $row.L
    $ws.Cells.Item($targetRow, 13).Value2 = $row.M
    $ws.Cells.Item($targetRow, 14).Value2 = $row.N
    $ws.Cells.Item($targetRow, 15).Value2 = $row.O
    $ws.Cells.Item($targetRow, 16).Value2 = $row.P
    $ws.Cells.Item($targetRow, 17).Value2 = $row.Q
    $ws.Cells.Item($targetRow, 18).Value2 = $row.R
    $ws.Cells.Item($targetRow, 19).Value2 = $row.S
    $ws.Cells.Item($targetRow, 20).Value2 = $row.T
    $ws.Cells.Item($targetRow, 21).Value2 = $row.U
    $ws.Cells.Item($targetRow, 22).Value2 = $row.V
}
